$wb = $excel.ActiveWorkbook

# Sheet "Overview" - Status and Latest Handoff Date columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-23-11 16:23:59"

# Sheet "zh-cn" - Status and Latest Handoff Datetime columns
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-11 16:23:56"

# Sheet "de-de" - Status and Latest Handoff Datetime columns
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-11 16:23:59"
